# Edit derived from the "Atualizacao de bases das ligas" commit diff.
# The source export re-ordered four pairs of adjacent match rows
# (rows 5/6, 16/17, 20/21, 69/70): every column except id (A),
# Div / Div Original Name (C/D) and Date (E) swaps between the two
# rows in each pair. (The shared-string-table side effects this causes
# - e.g. the lone G26 index change - do not change any resolved cell
# value, so they need no explicit action here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 5 <-> row 6
$ws.Range("B5").Value = 6781300
$ws.Range("B6").Value = 6781301
$ws.Range("F5").Value = "SV Altldersdorf"
$ws.Range("F6").Value = "FSV Saxonia Tangermunde"
$ws.Range("G5").Value = "SV Frankonia Wernsdorf"
$ws.Range("G6").Value = "MSC Preussen 1899"
$ws.Range("H5").Value = 8
$ws.Range("H6").Value = 4
$ws.Range("I5").Value = 2
$ws.Range("I6").Value = 0
$ws.Range("J5").Value = "H"
$ws.Range("J6").Value = "H"
$ws.Range("K5").Value = 2.2
$ws.Range("K6").Value = 1.65
$ws.Range("L5").Value = 3.5
$ws.Range("L6").Value = 4
$ws.Range("M5").Value = 2.7
$ws.Range("M6").Value = 4
$ws.Range("N5").Value = 1.727
$ws.Range("N6").Value = 1.45
$ws.Range("O5").Value = 4
$ws.Range("O6").Value = 4.5
$ws.Range("P5").Value = 3.5
$ws.Range("P6").Value = 5
$ws.Range("Q5").Value = -0.5
$ws.Range("Q6").Value = -1.25
$ws.Range("R5").Value = 1.775
$ws.Range("R6").Value = 1.925
$ws.Range("S5").Value = 2.025
$ws.Range("S6").Value = 1.875
$ws.Range("T5").Value = 3.25
$ws.Range("T6").Value = 3.5
$ws.Range("U5").Value = 1.925
$ws.Range("U6").Value = 1.825
$ws.Range("V5").Value = 1.875
$ws.Range("V6").Value = 1.975
$ws.Range("W5").Value = 0.7270000000000001
$ws.Range("W6").Value = 0.45
$ws.Range("X5").Value = -1
$ws.Range("X6").Value = -1
$ws.Range("Y5").Value = -1
$ws.Range("Y6").Value = -1
$ws.Range("Z5").Value = 0.7749999999999999
$ws.Range("Z6").Value = 0.925
$ws.Range("AA5").Value = -1
$ws.Range("AA6").Value = -1
$ws.Range("AB5").Value = 0.925
$ws.Range("AB6").Value = 0.825
$ws.Range("AC5").Value = -1
$ws.Range("AC6").Value = -1

# row 16 <-> row 17
$ws.Range("B16").Value = 7138608
$ws.Range("B17").Value = 7138607
$ws.Range("F16").Value = "SV UnterFlockenbach"
$ws.Range("F17").Value = "Rot Weiss Walldorf II"
$ws.Range("G16").Value = "SC Dortelweil"
$ws.Range("G17").Value = "Turnerschaft OberRoden"
$ws.Range("H16").Value = 1
$ws.Range("H17").Value = 3
$ws.Range("I16").Value = 1
$ws.Range("I17").Value = 2
$ws.Range("J16").Value = "D"
$ws.Range("J17").Value = "H"
$ws.Range("K16").Value = 1.083
$ws.Range("K17").Value = 2.25
$ws.Range("L16").Value = 9
$ws.Range("L17").Value = 3.75
$ws.Range("M16").Value = 16
$ws.Range("M17").Value = 2.5
$ws.Range("N16").Value = 1.125
$ws.Range("N17").Value = 2.25
$ws.Range("O16").Value = 7.5
$ws.Range("O17").Value = 3.8
$ws.Range("P16").Value = 13
$ws.Range("P17").Value = 2.45
$ws.Range("Q16").Value = -2.5
$ws.Range("Q17").Value = 0
$ws.Range("R16").Value = 1.775
$ws.Range("R17").Value = 1.8
$ws.Range("S16").Value = 1.925
$ws.Range("S17").Value = 2
$ws.Range("T16").Value = 4.25
$ws.Range("T17").Value = 3.75
$ws.Range("U16").Value = 1.975
$ws.Range("U17").Value = 1.95
$ws.Range("V16").Value = 1.825
$ws.Range("V17").Value = 1.85
$ws.Range("W16").Value = -1
$ws.Range("W17").Value = 1.25
$ws.Range("X16").Value = 6.5
$ws.Range("X17").Value = -1
$ws.Range("Y16").Value = -1
$ws.Range("Y17").Value = -1
$ws.Range("Z16").Value = -1
$ws.Range("Z17").Value = 0.8
$ws.Range("AA16").Value = 0.925
$ws.Range("AA17").Value = -1
$ws.Range("AB16").Value = -1
$ws.Range("AB17").Value = 0.95
$ws.Range("AC16").Value = 0.825
$ws.Range("AC17").Value = -1

# row 20 <-> row 21
$ws.Range("B20").Value = 7149361
$ws.Range("B21").Value = 7149166
$ws.Range("F20").Value = "TSG 1846 Bretzenheim"
$ws.Range("F21").Value = "Fuchse Berlin Reinickendorf"
$ws.Range("G20").Value = "TSV 1881 GauOdernheim"
$ws.Range("G21").Value = "SD Croatia Berlin"
$ws.Range("H20").Value = 4
$ws.Range("H21").Value = 1
$ws.Range("I20").Value = 2
$ws.Range("I21").Value = 4
$ws.Range("J20").Value = "H"
$ws.Range("J21").Value = "A"
$ws.Range("K20").Value = 3.75
$ws.Range("K21").Value = 1.4
$ws.Range("L20").Value = 4.333
$ws.Range("L21").Value = 4.8
$ws.Range("M20").Value = 1.615
$ws.Range("M21").Value = 5.25
$ws.Range("N20").Value = 3.75
$ws.Range("N21").Value = 1.4
$ws.Range("O20").Value = 4.333
$ws.Range("O21").Value = 5
$ws.Range("P20").Value = 1.615
$ws.Range("P21").Value = 5
$ws.Range("Q20").Value = 0.75
$ws.Range("Q21").Value = -1.25
$ws.Range("R20").Value = 2
$ws.Range("R21").Value = 1.8
$ws.Range("S20").Value = 1.8
$ws.Range("S21").Value = 2
$ws.Range("T20").Value = 3.75
$ws.Range("T21").Value = 3.5
$ws.Range("U20").Value = 1.9
$ws.Range("U21").Value = 1.925
$ws.Range("V20").Value = 1.9
$ws.Range("V21").Value = 1.875
$ws.Range("W20").Value = 2.75
$ws.Range("W21").Value = -1
$ws.Range("X20").Value = -1
$ws.Range("X21").Value = -1
$ws.Range("Y20").Value = -1
$ws.Range("Y21").Value = 4
$ws.Range("Z20").Value = 1
$ws.Range("Z21").Value = -1
$ws.Range("AA20").Value = -1
$ws.Range("AA21").Value = 1
$ws.Range("AB20").Value = 0.8999999999999999
$ws.Range("AB21").Value = 0.925
$ws.Range("AC20").Value = -1
$ws.Range("AC21").Value = -1

# row 69 <-> row 70
$ws.Range("B69").Value = 7423699
$ws.Range("B70").Value = 7423700
$ws.Range("F69").Value = "SG 2000 MulheimKarlich"
$ws.Range("F70").Value = "TuS Hornau"
$ws.Range("G69").Value = "Ahrweiler BC"
$ws.Range("G70").Value = "FC Burgsolms"
$ws.Range("H69").Value = 2
$ws.Range("H70").Value = 3
$ws.Range("I69").Value = 2
$ws.Range("I70").Value = 0
$ws.Range("J69").Value = "D"
$ws.Range("J70").Value = "H"
$ws.Range("K69").Value = 2.2
$ws.Range("K70").Value = 1.727
$ws.Range("L69").Value = 5
$ws.Range("L70").Value = 4.5
$ws.Range("M69").Value = 2.2
$ws.Range("M70").Value = 3.2
$ws.Range("N69").Value = 2.2
$ws.Range("N70").Value = 1.727
$ws.Range("O69").Value = 4.75
$ws.Range("O70").Value = 4.5
$ws.Range("P69").Value = 2.2
$ws.Range("P70").Value = 3.2
$ws.Range("Q69").Value = 0
$ws.Range("Q70").Value = -0.5
$ws.Range("R69").Value = 1.9
$ws.Range("R70").Value = 1.775
$ws.Range("S69").Value = 1.9
$ws.Range("S70").Value = 2.025
$ws.Range("T69").Value = 4.25
$ws.Range("T70").Value = 3.5
$ws.Range("U69").Value = 1.775
$ws.Range("U70").Value = 1.85
$ws.Range("V69").Value = 2.025
$ws.Range("V70").Value = 1.95
$ws.Range("W69").Value = -1
$ws.Range("W70").Value = 0.7270000000000001
$ws.Range("X69").Value = 3.75
$ws.Range("X70").Value = -1
$ws.Range("Y69").Value = -1
$ws.Range("Y70").Value = -1
$ws.Range("Z69").Value = 0
$ws.Range("Z70").Value = 0.7749999999999999
$ws.Range("AA69").Value = -0
$ws.Range("AA70").Value = -1
$ws.Range("AB69").Value = -0.5
$ws.Range("AB70").Value = -1
$ws.Range("AC69").Value = 0.5125
$ws.Range("AC70").Value = 0.95

